# Apply the data update described by the commit "make a little change".
# Net effect: the first two existing rows (id 15, id 16) are gone, the
# remaining two existing rows (id 17, id 18) shift up to rows 2-3 (id 17's
# row also gets a few field updates), and two brand-new rows (id 19, id 20)
# are appended at rows 4-5. The header row (row 1) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was id 17 / "Naomi" in row 4) -> now id 17 / "Naomi clara" with a
# few updated fields (product code + total).
$ws.Range("A2").Value = 17
$ws.Range("B2").Value = "Naomi clara"
$ws.Range("C2").Value = "Jl. merdeka no 17"
$ws.Range("D2").Value = "Jawa Tengah"
$ws.Range("E2").Value = "Semarang"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "08963764826"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "2022-11-10 20:26"
$ws.Range("H2").Value = "SG005LTR"
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 4550000

# Row 3 (was id 18 / "Lestari" in row 5) -> unchanged values, just shifted up.
$ws.Range("A3").Value = 18
$ws.Range("B3").Value = "Lestari"
$ws.Range("C3").Value = "cilandak"
$ws.Range("D3").Value = "Daerah Khusus Ibukota Jakarta"
$ws.Range("E3").Value = "Jakarta"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "39480579170"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "2022-11-10 00:00:00"
$ws.Range("H3").Value = "FL007KG"
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 6720000

# Row 4 - brand-new record, id 19 / "Nuvo".
$ws.Range("A4").Value = 19
$ws.Range("B4").Value = "Nuvo"
$ws.Range("C4").Value = "faldhiuagiuga"
$ws.Range("D4").Value = "Daerah Khusus Ibukota Jakarta"
$ws.Range("E4").Value = "Jakarta"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "08964384729"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "2022-11-11 00:00:00"
$ws.Range("H4").Value = "FL007KG"
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 22400000

# Row 5 - brand-new record, id 20 / "Jeremi".
$ws.Range("A5").Value = 20
$ws.Range("B5").Value = "Jeremi"
$ws.Range("C5").Value = "Balik Papan"
$ws.Range("D5").Value = "Kalimantan Selatan"
$ws.Range("E5").Value = "Balik Papan"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "0897463532"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "2022-11-11 00:00:00"
$ws.Range("H5").Value = "MFHF1C1"
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2130000
